# New order came in at 2026-01-13 16:54 (order #5, "Pooja" / Girl Holding
# Hands Thali x1). Insert it as the new top data row on "All Orders" (row 2),
# which pushes the previously-existing rows 2-5 down to 3-6, and bump the
# "Total Orders" count on the "Daily Summary" sheet from 4 to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Push the existing data down one row and open up a fresh row 2.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "2026-01-13 16:54"
$ws.Range("C2").Value = "Pooja"

# The phone number is purely numeric text ("9096648553") and must stay a
# text value (matching the other Phone cells in this column) instead of
# being auto-converted to a number. Force text via NumberFormat, write the
# value, then drop the now-unneeded number format so the cell is left with
# the workbook's default (unformatted) style, same as its neighbours.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "9096648553"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "Level 1, Tower S3, CyberCity,`nMagarpatta City, Hadapsar, PUNE 411013"
$ws.Range("F2").Value = "Girl Holding Hands Thali x1"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# Reflect the new order in the daily roll-up: Total Orders 4 -> 5.
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Range("B2").Value = 5
